# Generate Report for Handoff
#
# Refreshes the localization-status report after a new handoff run:
# the four "Ready for handoff" source files (rows 4-7) in each language
# sheet move from priority "low" to "ht", and their handoff timestamp is
# bumped to the time the new .xlf files were generated. The Overview
# sheet's "Latest HO Xliff Generate Date" column mirrors the de-de sheet's
# handoff timestamp for those same rows, so it is refreshed too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$zhCnHandoffDatetime = "2016-09-07 02:38:56"
$deDeHandoffDatetime = "2016-09-07 02:39:07"

for ($row = 4; $row -le 7; $row++) {
    # Priority: low -> ht
    $wsZhCn.Cells.Item($row, 5).Value = "ht"
    $wsDeDe.Cells.Item($row, 5).Value = "ht"

    # Latest Handoff Datetime (column H) refreshed for the new handoff
    $wsZhCn.Cells.Item($row, 8).Value = $zhCnHandoffDatetime
    $wsDeDe.Cells.Item($row, 8).Value = $deDeHandoffDatetime

    # Overview sheet's Latest HO Xliff Generate Date (column G) mirrors
    # the de-de handoff datetime for these rows
    $wsOverview.Cells.Item($row, 7).Value = $deDeHandoffDatetime
}
